# Add "Baptismal Date" column to the confirmation detail sheet.
#
# The new column is inserted between "Birth Date" (B) and "Father" (old C),
# shifting Father/Mother/Sponsor 1/Sponsor 2/Contact Number/Present Address
# one column to the right (C->D, D->E, E->F, F->G, G->H, H->I).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at C; this shifts the old C..J columns to D..K and
# keeps each existing cell's formatting (it copies format from the column
# on the left, same as Excel's default "Insert Sheet Columns" behaviour).
$ws.Columns.Item(3).Insert()

# Give the new column the same cell formatting (date number format, centered
# header style, etc.) as column B (Birth Date) by copying B1:B5 -> C1:C5.
$ws.Range("B1:B5").Copy()
$ws.Range("C1:C5").PasteSpecial(-4122)

# Match the new column's width to column B's width as closely as the
# Excel object model allows.
$ws.Columns.Item(3).ColumnWidth = $ws.Columns.Item(2).ColumnWidth()

# Header text for the new column.
$ws.Cells.Item(2, 3).Value2 = "Baptismal Date"

# Baptismal date values for the three detail rows (serial dates).
$ws.Cells.Item(3, 3).Value2 = 45352
$ws.Cells.Item(4, 3).Value2 = 45353
$ws.Cells.Item(5, 3).Value2 = 45354

# Restore the cursor/selection position recorded in the saved file.
$ws.Range("F11").Select()

Write-Output "Added Baptismal Date column"
